# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, rows 2 and 5-13.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 6547
    5  = 48
    6  = 1969
    7  = 1503
    8  = 305
    9  = 1004
    10 = 380
    11 = 7
    12 = 5625
    13 = 75
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
